$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("4:4").Insert()

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C4").Value = "Los Lagos"
$ws.Range("D4").Value = 45190
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 100112012
$ws.Range("G4").Value = "Espinaca"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 12000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 12000
$ws.Range("N4").Value = "$/cuna 10 kilos"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 1200
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = "Hortaliza"
